$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so Excel keeps the exact
# literal string (e.g. "577.06") instead of coercing to a float, matching
# the inline-string values recorded in the source feed dump.
$textCells = @("D5","D6","D12","D16","D19","D20","D21","D22","D23","D24","D29","D32","D34","D37","D38","D41","D42","D43","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.637.07'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '3.229.56'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '577.06'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '175.10'
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '3.227.36'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '0.391'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '3.797.60'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('E14').Value = '  -2.94%  '
$ws.Range('D15').Value = '64.835.43'
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '25.69'
$ws.Range('E16').Value = '  -3.12%  '
$ws.Range('D17').Value = '3.228.99'
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('E18').Value = '  -2.30%  '
$ws.Range('D19').Value = '415.32'
$ws.Range('E19').Value = '  -4.02%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '12.85'
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = '5.37'
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').Value = '7.20'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '70.34'
$ws.Range('E24').Value = '  -2.24%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  +4.08%  '
$ws.Range('E27').Value = '  -1.99%  '
$ws.Range('E28').Value = '  -2.43%  '
$ws.Range('D29').Value = '9.15'
$ws.Range('E29').Value = '  +2.94%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').Value = '21.85'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '5.01'
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('D37').Value = '156.13'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '1.39'
$ws.Range('E38').Value = '  -2.68%  '
$ws.Range('D39').Value = '2.829.25'
$ws.Range('E39').Value = '  +2.13%  '
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('D41').Value = '25.46'
$ws.Range('E41').Value = '  -4.31%  '
$ws.Range('D42').Value = '4.21'
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('D43').Value = '0.728'
$ws.Range('E43').Value = '  -6.19%  '
$ws.Range('E44').Value = '  -2.68%  '
$ws.Range('E45').Value = '  -4.67%  '
$ws.Range('D46').Value = '0.0626'
$ws.Range('E46').Value = '  -4.61%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.18'
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '303.86'
$ws.Range('E48').Value = '  -5.57%  '
$ws.Range('D49').Value = '22.25'
$ws.Range('E49').Value = '  -4.45%  '
$ws.Range('D50').Value = '0.0263'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('E51').Value = '  -0.52%  '
